# Add a new "Flowering_initiation_exp2.csv" column-description section
# to the Metadata sheet, right after the existing "Flower_size_exp2.csv"
# section (which ends at row 85).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Section header (bold, like the other "Column descriptions for ..." rows)
$ws.Range("A87").Value = 'Column descriptions for "Flowering_initiation_exp2.csv" '
$ws.Range("A87").Font.Bold = $true

# Column description rows (row 88-93)
$ws.Range("A88").Value = "Date"
$ws.Range("B88").Value = "Date that flowering initiated"

$ws.Range("A89").Value = "Chamber"
$ws.Range("B89").Value = "Greenhouse chamber identity (60, 61, 62, 63, or 65)"

$ws.Range("A90").Value = "CO2"
$ws.Range("B90").Value = "CO2 treatment (elevated or ambient)"

$ws.Range("A91").Value = "Round"
$ws.Range("B91").Value = "Experimental round (1 or 2)"

$ws.Range("A92").Value = "Plant"
$ws.Range("B92").Value = "Plant species common name abbreviation (B = borage, BW = buckwheat, C = red clover, N = nasturtium, LP = lacy phacelia, PP, = partridge pea, SF = sunflower, SA = sweet alyssum)"

$ws.Range("A93").Value = "Plant_ID"
$ws.Range("B93").Value = "Unique sample identifier for the individual plant where the sample came from"

# Update view selection to reflect the new last-edited cell (B89), matching
# the document's recorded cursor position after the edit.
$null = $ws.Range("B89").Select()
